# Weekly Fruta/Hortaliza update: insert a new price record for Damasco
# (Castle Brite, Región de O'Higgins) as row 18, pushing the existing
# rows 18-21 down to 19-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 18 (existing rows 18-21 shift to 19-22)
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with the new record
$ws.Cells.Item(18, 1).Value = 7
$ws.Cells.Item(18, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(18, 3).Value = "Ñuble"
$ws.Cells.Item(18, 4).Value = 45275
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103003
$ws.Cells.Item(18, 10).Value = "Damasco"
$ws.Cells.Item(18, 11).Value = "Castle Brite"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 80
$ws.Cells.Item(18, 14).Value = 17000
$ws.Cells.Item(18, 15).Value = 18000
$ws.Cells.Item(18, 16).Value = 17375
$ws.Cells.Item(18, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(18, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 19).Value = 1738
$ws.Cells.Item(18, 20).Value = 10
